# "Add files via upload" — [D] results of finding edge of chaos
# Re-samples the "edge of chaos" region (just above the previous coarse
# step) on both the "random" and "periodic" sheets, inserting finer-grained
# x values (column A) with their corresponding measurements (column F)
# around the transition, and nudges the active sheet / selection / scroll
# position to match the author's last-saved view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "random" (sheet1) — rows 15-32 get new, finer-grained samples.
# ---------------------------------------------------------------------
$wsRandom = $wb.Worksheets.Item("random")

# Clear out the old F-column measurements that fell in rows 18-22; the
# new layout leaves those rows (57.3 .. 57.7) without a measurement.
$wsRandom.Range("F18:F22").ClearContents()

$randomData = @(
    @(15, 57,   -1.1),
    @(16, 57.1, -0.94),
    @(17, 57.2, -0.88),
    @(18, 57.3, $null),
    @(19, 57.4, $null),
    @(20, 57.5, $null),
    @(21, 57.6, $null),
    @(22, 57.7, $null),
    @(23, 57.8, $null),
    @(24, 57.9, $null),
    @(25, 58,   $null),
    @(26, 59,   $null),
    @(27, 59.1, 1.15),
    @(28, 60,   1.65),
    @(29, 70,   4.22),
    @(30, 80,   8.17),
    @(31, 90,   8.58),
    @(32, 100,  8.27)
)

foreach ($entry in $randomData) {
    $r = $entry[0]
    $wsRandom.Cells.Item($r, 1).Value = $entry[1]
    if ($entry[2] -ne $null) {
        $wsRandom.Cells.Item($r, 6).Value = $entry[2]
    }
}

# Row 15's new measurement is highlighted in red, like the other flagged
# edge-of-chaos points elsewhere in the workbook.
$wsRandom.Range("F15").Font.Color = 255

# ---------------------------------------------------------------------
# Sheet "periodic" (sheet2) — rows 11-27 get new, finer-grained samples.
# ---------------------------------------------------------------------
$wsPeriodic = $wb.Worksheets.Item("periodic")

# Row 13 loses its old measurement; the new row 13 (52.3) has none.
$wsPeriodic.Range("F13").ClearContents()

$periodicData = @(
    @(11, 52.1, -0.41),
    @(12, 52.2, -0.47),
    @(13, 52.3, $null),
    @(14, 52.4, $null),
    @(15, 52.5, $null),
    @(16, 52.6, $null),
    @(17, 52.7, $null),
    @(18, 52.8, $null),
    @(19, 52.9, $null),
    @(20, 53,   0.45),
    @(21, 54,   0.46),
    @(22, 55,   -2),
    @(23, 56,   $null),
    @(24, 57,   $null),
    @(25, 58,   $null),
    @(26, 59,   $null),
    @(27, 60,   1.4)
)

foreach ($entry in $periodicData) {
    $r = $entry[0]
    $wsPeriodic.Cells.Item($r, 1).Value = $entry[1]
    if ($entry[2] -ne $null) {
        $wsPeriodic.Cells.Item($r, 6).Value = $entry[2]
    }
}

# ---------------------------------------------------------------------
# View state: "random" becomes the active/selected tab (it was
# "periodic" before), with the selection/scroll position the author
# left things in on each sheet.
# ---------------------------------------------------------------------
$wsRandom.Activate()
$wsRandom.Range("H26").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$wsPeriodic.Range("I15").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wsRandom.Activate()
